# Amend corrected label annotations
# The F column ("labels") contains human-entered category labels that were
# normalized to lowercase. One cell (row 32) additionally had its two
# "||"-separated segments reordered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)   # column F
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -eq "") { continue }

    if ($r -eq 32) {
        $cell.Value2 = "application instructions || env warning - species"
    } else {
        $cell.Value2 = $val.ToLower()
    }
}
